# daily auto push: 2026-01-30 07:04 UTC
# Insert a new data row for 2026/01/30 07:00 bucket right before the
# existing 2026/12/29 13:00 row (row 749), shifting all subsequent rows
# down by one, then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 749:790 down to 750:791 by inserting a new row at 749.
$ws.Rows.Item(749).Insert()

# Populate the newly inserted row. The leading apostrophe forces the
# date-like text to be stored as text (matching the rest of column A)
# instead of being auto-converted to a date serial number; resetting
# the style back to Normal afterward removes the quote-prefix marker
# so the cell formatting matches its neighbors exactly.
$ws.Range("A749").Value2 = "'2026/01/30"
$ws.Range("A749").Style = "Normal"
$ws.Range("B749").Value2 = "金"
$ws.Range("C749").Value2 = 13
$ws.Range("D749").Value2 = 190
